# Weekly update: insert a new price record as row 116, shifting the
# existing rows (116-225) down by one (to 117-226).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 116 (pushes rows 116:225 down to 117:226,
# and keeps column formatting from the row above, matching Excel's
# default "Insert" behaviour).
$ws.Rows("116:116").Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(116, 1).Value = 3
$ws.Cells.Item(116, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44589
$ws.Cells.Item(116, 5).Value = 5
$ws.Cells.Item(116, 6).Value = 100112001
$ws.Cells.Item(116, 7).Value = "Berenjena"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 73
$ws.Cells.Item(116, 11).Value = 9500
$ws.Cells.Item(116, 12).Value = 10000
$ws.Cells.Item(116, 13).Value = 9740
$ws.Cells.Item(116, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 162
$ws.Cells.Item(116, 17).Value = 60
$ws.Cells.Item(116, 18).Value = "Hortaliza"
